# Apply the edit described by the commit "table ss with pre-experiment data":
#  - relabel the ">100000" row header (used for the pre-experiment subsample)
#    as LaTeX math "$>$100000"
#  - add a thin top border above that row to visually separate the new block
#  - restore the active-cell/selection anchor to A2 (top-left) over A2:E16

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the ">100000" labels in row 10 to the LaTeX-escaped "$>$100000"
$ws.Range("B10:E10").Value = '$>$100000'

# 2) Add a thin top border across A10:E10 to set the row apart
$topRow = $ws.Range("A10:E10")
$topRow.Borders.Item(8).LineStyle = 1
$topRow.Borders.Item(8).Weight = 2

# 3) Reset the sheet selection to match the saved view (active cell A2, same range selected)
$ws.Range("A2").Activate() | Out-Null
$ws.Range("A2:E16").Select() | Out-Null
